$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 3 (E3 text change, F3/G3 new values)
$ws.Range("E3").Value = "`$.get('/items/:id')"
$ws.Range("F3").Value = "`$.ajax({})"
$ws.Range("G3").Value = "fetch('/items/:id')"

# Row 4: render a form / JavaScript
$ws.Range("A4").Value = "render a form"
$ws.Range("C4").Value = "JavaScript"

# Row 5: render a form / HTML from Rails
$ws.Range("A5").Value = "render a form"
$ws.Range("C5").Value = "HTML from Rails"

# Row 6: render a form / Handlebars
$ws.Range("A6").Value = "render a form"
$ws.Range("C6").Value = "Handlebars"

# Column width adjustments
$ws.Columns.Item(3).ColumnWidth = 24.5
$ws.Columns.Item(7).ColumnWidth = 26

# Selection change
$ws.Range("H5").Select()
